$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "35.448.21"
$ws.Range("E2").Value = "  +0.65%  "

$ws.Range("D3").Value = "1.923.67"
$ws.Range("E3").Value = "  +1.68%  "

$ws.Range("E4").Value = "  +0.09%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "0.731"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +11.65%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "253.89"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +4.62%  "

$ws.Range("E7").Value = "  +0.09%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "40.97"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -0.68%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.357"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +2.60%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "52.62"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +4.52%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0748"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +5.79%  "

$ws.Range("E12").Value = "  +0.50%  "

$ws.Range("D13").Value = "2.206.89"
$ws.Range("E13").Value = "  +1.76%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "12.82"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +7.47%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.720"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +4.31%  "

$ws.Range("D16").Value = "1.925.77"
$ws.Range("E16").Value = "  +1.82%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "4.91"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +1.94%  "

$ws.Range("D18").Value = "35.485.85"
$ws.Range("E18").Value = "  +0.74%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "74.46"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +4.73%  "

$ws.Range("D20").Value = "0.0₃0837"
$ws.Range("E20").Value = "  +3.37%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "243.64"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +1.28%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "13.04"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +5.48%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "5.09"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +7.96%  "

$ws.Range("E24").Value = "  -0.02%  "

$ws.Range("B25").Value = "Toncoin"
$ws.Range("C25").Value = "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.34"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +2.62%  "

$ws.Range("B26").Value = "PancakeSwap"
$ws.Range("C26").Value = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "2.41"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -0.34%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "168.04"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -1.22%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "8.65"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +3.39%  "

$ws.Range("E29").Value = "  +6.62%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "18.80"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +3.43%  "

$ws.Range("D31").Value = "4.128.20"
$ws.Range("E31").Value = "  +19.44%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "4.38"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +6.91%  "

$ws.Range("B33").Value = "WEMIXToken"
$ws.Range("C33").Value = "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "1.98"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +14.35%  "

$ws.Range("B34").Value = "TrustWalletToken"
$ws.Range("C34").Value = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.63"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +23.94%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.0583"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +4.35%  "

$ws.Range("E36").Value = "  +3.81%  "

$ws.Range("E37").Value = "  -0.01%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.918"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -2.01%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "2.04"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +1.33%  "

$ws.Range("E40").Value = "  +10.29%  "

$ws.Range("B41").Value = "Aave"
$ws.Range("C41").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "97.43"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +9.66%  "

$ws.Range("B42").Value = "ARBITRUM"
$ws.Range("C42").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "1.13"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +4.24%  "

$ws.Range("E43").Value = "  +1.79%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.0654"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +2.64%  "

$ws.Range("D45").Value = "1.349.10"
$ws.Range("E45").Value = "  +1.00%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "2.46"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +4.85%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "6.86"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +5.61%  "

$ws.Range("E48").Value = "  +0.73%  "

$ws.Range("E49").Value = "  +0.64%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "45.25"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -6.35%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "11.92"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +5.20%  "
